$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.93209999999999
$ws.Range("A9").Value = -20.46069999999997
$ws.Range("B11").Value = 5.366599999999995
$ws.Range("A18").Value = -23.00020000000001
$ws.Range("A20").Value = -22.19010000000003
$ws.Range("E21").Value = 13.0721
